$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.086674137469346
$ws.Range("D2").Value = 1.089176986687104
$ws.Range("E2").Value = 1.089112650765332
$ws.Range("F2").Value = 1.10028235738999
$ws.Range("I2").Value = 1.07130208120365
$ws.Range("J2").Value = 1.09152294714848
$ws.Range("K2").Value = 1.091827776155915
$ws.Range("L2").Value = 1.091763605515123
$ws.Range("M2").Value = 1.102904928915404
$ws.Range("N2").Value = 1.093073034932863
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.088138859816106
$ws.Range("D3").Value = 1.090206265893598
$ws.Range("E3").Value = 1.090392534083969
$ws.Range("F3").Value = 1.101525422823523
$ws.Range("I3").Value = 1.071829409084944
$ws.Range("J3").Value = 1.092648973132491
$ws.Range("K3").Value = 1.092676436901878
$ws.Range("L3").Value = 1.092862261228252
$ws.Range("M3").Value = 1.103968913522419
$ws.Range("N3").Value = 1.094200660002928
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.089085242288455
$ws.Range("D4").Value = 1.090871291718197
$ws.Range("E4").Value = 1.091219195359617
$ws.Range("F4").Value = 1.102328715229903
$ws.Range("I4").Value = 1.072168729291139
$ws.Range("J4").Value = 1.093375705332727
$ws.Range("K4").Value = 1.093223978314146
$ws.Range("L4").Value = 1.093571094736764
$ws.Range("M4").Value = 1.104655756013036
$ws.Range("N4").Value = 1.094928424246246
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.089482775171832
$ws.Range("D5").Value = 1.091150635698331
$ws.Range("E5").Value = 1.091566367948941
$ws.Range("F5").Value = 1.102666171606144
$ws.Range("I5").Value = 1.072310928051395
$ws.Range("J5").Value = 1.093680778117489
$ws.Range("K5").Value = 1.093453785250114
$ws.Range("L5").Value = 1.093868597784011
$ws.Range("M5").Value = 1.104944119503632
$ws.Range("N5").Value = 1.095233930269355
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.089549503750289
$ws.Range("D6").Value = 1.09119752522212
$ws.Range("E6").Value = 1.091624639045918
$ws.Range("F6").Value = 1.102722817621281
$ws.Range("I6").Value = 1.07233477744779
$ws.Range("J6").Value = 1.09373197520782
$ws.Range("K6").Value = 1.093492348688628
$ws.Range("L6").Value = 1.093918521234404
$ws.Range("M6").Value = 1.104992514546312
$ws.Range("N6").Value = 1.095285200065426
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.089090555416198
$ws.Range("D7").Value = 1.090875025237055
$ws.Range("E7").Value = 1.09122383568813
$ws.Range("F7").Value = 1.102333225306316
$ws.Range("I7").Value = 1.07217063112875
$ws.Range("J7").Value = 1.093379783474128
$ws.Range("K7").Value = 1.093227050492669
$ws.Range("L7").Value = 1.093575071908338
$ws.Range("M7").Value = 1.104659610647877
$ws.Range("N7").Value = 1.094932508179076
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.087169439241698
$ws.Range("D8").Value = 1.089525042716273
$ws.Range("E8").Value = 1.089545508994719
$ws.Range("F8").Value = 1.100702678191781
$ws.Range("I8").Value = 1.071480688257489
$ws.Range("J8").Value = 1.091903885640143
$ws.Range("K8").Value = 1.092114918065578
$ws.Range("L8").Value = 1.092135333053762
$ws.Range("M8").Value = 1.10326484693358
$ws.Range("N8").Value = 1.093454514400878
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.083773235947076
$ws.Range("D9").Value = 1.087138494432813
$ws.Range("E9").Value = 1.08657628443614
$ws.Range("F9").Value = 1.09782116657297
$ws.Range("I9").Value = 1.070250285764345
$ws.Range("J9").Value = 1.089288522725351
$ws.Range("K9").Value = 1.090142805227712
$ws.Range("L9").Value = 1.089582246996677
$ws.Range("M9").Value = 1.100794441212153
$ws.Range("N9").Value = 1.090835437370761
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.081501317524365
$ws.Range("D10").Value = 1.085542068436278
$ws.Range("E10").Value = 1.084588527265302
$ws.Range("F10").Value = 1.095894301205078
$ws.Range("I10").Value = 1.06942001914685
$ws.Range("J10").Value = 1.087534780377982
$ws.Range("K10").Value = 1.088819521427007
$ws.Range("L10").Value = 1.087869057646891
$ws.Range("M10").Value = 1.099138721152878
$ws.Range("N10").Value = 1.089079204508102
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.080515611773591
$ws.Range("D11").Value = 1.084849471005537
$ws.Range("E11").Value = 1.083725767345164
$ws.Range("F11").Value = 1.095058492387221
$ws.Range("I11").Value = 1.069058097373697
$ws.Range("J11").Value = 1.086772906776822
$ws.Range("K11").Value = 1.088244451709923
$ws.Range("L11").Value = 1.087124515801567
$ws.Range("M11").Value = 1.098419631026077
$ws.Range("N11").Value = 1.08831624895904
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.080149175613176
$ws.Range("D12").Value = 1.084592004768918
$ws.Range("E12").Value = 1.08340498552744
$ws.Range("F12").Value = 1.094747810185615
$ws.Range("I12").Value = 1.068923298166615
$ws.Range("J12").Value = 1.086489532349238
$ws.Range("K12").Value = 1.088030528699602
$ws.Range("L12").Value = 1.086847544972318
$ws.Range("M12").Value = 1.098152200139817
$ws.Range("N12").Value = 1.088032472107265
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.080227791173326
$ws.Range("D13").Value = 1.084647241482626
$ws.Range("E13").Value = 1.083473808649998
$ws.Range("F13").Value = 1.094814462873532
$ws.Range("I13").Value = 1.068952229637502
$ws.Range("J13").Value = 1.086550334452972
$ws.Range("K13").Value = 1.088076430301345
$ws.Range("L13").Value = 1.086906975028534
$ws.Range("M13").Value = 1.098209579938872
$ws.Range("N13").Value = 1.088093360556961
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.080485328237253
$ws.Range("D14").Value = 1.084828192966767
$ws.Range("E14").Value = 1.083699257871922
$ws.Range("F14").Value = 1.095032815934846
$ws.Range("I14").Value = 1.069046962300397
$ws.Range("J14").Value = 1.086749490760067
$ws.Range("K14").Value = 1.08822677524804
$ws.Range("L14").Value = 1.087101629805304
$ws.Range("M14").Value = 1.098397531844414
$ws.Range("N14").Value = 1.088292799688856
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.080643965170454
$ws.Range("D15").Value = 1.084939655894752
$ws.Range("E15").Value = 1.083838122720524
$ws.Range("F15").Value = 1.095167320397087
$ws.Range("I15").Value = 1.069105281714283
$ws.Range("J15").Value = 1.08687214690439
$ws.Range("K15").Value = 1.08831936566601
$ws.Range("L15").Value = 1.087221507879901
$ws.Range("M15").Value = 1.0985132914852
$ws.Range("N15").Value = 1.088415630018972
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.081566693785903
$ws.Range("D16").Value = 1.085588005371514
$ws.Range("E16").Value = 1.084645742082047
$ws.Range("F16").Value = 1.095949739763135
$ws.Range("I16").Value = 1.069443987641381
$ws.Range("J16").Value = 1.087585290386613
$ws.Range("K16").Value = 1.088857642752351
$ws.Range("L16").Value = 1.087918412606279
$ws.Range("M16").Value = 1.099186398985712
$ws.Range("N16").Value = 1.089129786246738
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.082144969218944
$ws.Range("D17").Value = 1.085994337627199
$ws.Range("E17").Value = 1.085151787369005
$ws.Range("F17").Value = 1.09644013493227
$ws.Range("I17").Value = 1.069655801256036
$ws.Range("J17").Value = 1.088031954936254
$ws.Range("K17").Value = 1.089194730192652
$ws.Range("L17").Value = 1.0883548300246
$ws.Range("M17").Value = 1.099608041460622
$ws.Range("N17").Value = 1.089577085111273
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.082482079896657
$ws.Range("D18").Value = 1.086231215891535
$ws.Range("E18").Value = 1.085446757880787
$ws.Range("F18").Value = 1.096726033331127
$ws.Range("I18").Value = 1.069779116155787
$ws.Range("J18").Value = 1.088292246838096
$ws.Range("K18").Value = 1.089391147273303
$ws.Range("L18").Value = 1.088609122949098
$ws.Range("M18").Value = 1.099853771074547
$ws.Range("N18").Value = 1.089837746657471
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.082596994352345
$ws.Range("D19").Value = 1.086311963595227
$ws.Range("E19").Value = 1.085547301965136
$ws.Range("F19").Value = 1.096823493586116
$ws.Range("I19").Value = 1.069821124025688
$ws.Range("J19").Value = 1.08838095913665
$ws.Range("K19").Value = 1.089458086568988
$ws.Range("L19").Value = 1.088695786018177
$ws.Range("M19").Value = 1.099937523522708
$ws.Range("N19").Value = 1.089926584937665
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.082082945201268
$ws.Range("D20").Value = 1.085950755340205
$ws.Range("E20").Value = 1.085097513938795
$ws.Range("F20").Value = 1.096387534777585
$ws.Range("I20").Value = 1.069633099720614
$ws.Range("J20").Value = 1.087984056924637
$ws.Range("K20").Value = 1.089158584650681
$ws.Range("L20").Value = 1.088308033708218
$ws.Range("M20").Value = 1.099562824695642
$ws.Range("N20").Value = 1.089529119078985
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.080409498370204
$ws.Range("D21").Value = 1.08477491295
$ws.Range("E21").Value = 1.083632877427304
$ws.Range("F21").Value = 1.094968522674304
$ws.Range("I21").Value = 1.069019076003016
$ws.Range("J21").Value = 1.086690854750206
$ws.Range("K21").Value = 1.088182511179279
$ws.Range("L21").Value = 1.08704432032717
$ws.Range("M21").Value = 1.098342193812902
$ws.Range("N21").Value = 1.088234080409134
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.079355589442367
$ws.Range("D22").Value = 1.084034427644236
$ws.Range("E22").Value = 1.08271018013346
$ws.Range("F22").Value = 1.094075025366734
$ws.Range("I22").Value = 1.068630899569041
$ws.Range("J22").Value = 1.085875561734428
$ws.Range("K22").Value = 1.087566981724887
$ws.Range("L22").Value = 1.086247370319332
$ws.Range("M22").Value = 1.097572829706997
$ws.Range("N22").Value = 1.087417629583759
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.079914454796121
$ws.Range("D23").Value = 1.084427086804483
$ws.Range("E23").Value = 1.083199494393251
$ws.Range("F23").Value = 1.094548811110063
$ws.Range("I23").Value = 1.068836880806989
$ws.Range("J23").Value = 1.08630797516144
$ws.Range("K23").Value = 1.087893460706394
$ws.Range("L23").Value = 1.086670078496041
$ws.Range("M23").Value = 1.097980866569505
$ws.Range("N23").Value = 1.087850657087436
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.082110971765011
$ws.Range("D24").Value = 1.085970448696984
$ws.Range("E24").Value = 1.085122038375184
$ws.Range("F24").Value = 1.096411302958747
$ws.Range("I24").Value = 1.069643358285537
$ws.Range("J24").Value = 1.088005700716091
$ws.Range("K24").Value = 1.08917491788445
$ws.Range("L24").Value = 1.08832917975899
$ws.Range("M24").Value = 1.099583256846242
$ws.Range("N24").Value = 1.089550793607106
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.084652576750772
$ws.Range("D25").Value = 1.087756409392577
$ws.Range("E25").Value = 1.087345333361903
$ws.Range("F25").Value = 1.098567117376155
$ws.Range("I25").Value = 1.070570124878832
$ws.Range("J25").Value = 1.089966424985054
$ws.Range("K25").Value = 1.090654133798796
$ws.Range("L25").Value = 1.090244219506251
$ws.Range("M25").Value = 1.101434628066126
$ws.Range("N25").Value = 1.091514302329429
